# Update the answer table with the newly generated values.
# Each cell's Range.Text is set directly (by row/column) so there is no
# ambiguity from overlapping old/new text values during a Find/Replace pass.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "25÷5=5, 0" },
    @{ Row = 1;  Col = 2; Text = "79÷7=11, 2" },
    @{ Row = 1;  Col = 3; Text = "20÷8=2, 4" },
    @{ Row = 1;  Col = 4; Text = "77÷6=12, 5" },
    @{ Row = 1;  Col = 5; Text = "42÷8=5, 2" },

    @{ Row = 5;  Col = 1; Text = "63÷2=31, 1" },
    @{ Row = 5;  Col = 2; Text = "84÷5=16, 4" },
    @{ Row = 5;  Col = 3; Text = "34÷4=8, 2" },
    @{ Row = 5;  Col = 4; Text = "80÷7=11, 3" },
    @{ Row = 5;  Col = 5; Text = "64÷7=9, 1" },

    @{ Row = 9;  Col = 1; Text = "77÷9=8, 5" },
    @{ Row = 9;  Col = 2; Text = "11÷2=5, 1" },
    @{ Row = 9;  Col = 3; Text = "92÷9=10, 2" },
    @{ Row = 9;  Col = 4; Text = "30÷3=10, 0" },
    @{ Row = 9;  Col = 5; Text = "39÷5=7, 4" },

    @{ Row = 13; Col = 1; Text = "52÷2=26, 0" },
    @{ Row = 13; Col = 2; Text = "45÷2=22, 1" },
    @{ Row = 13; Col = 3; Text = "35÷4=8, 3" },
    @{ Row = 13; Col = 4; Text = "92÷5=18, 2" },
    @{ Row = 13; Col = 5; Text = "32÷5=6, 2" },

    @{ Row = 17; Col = 1; Text = "46÷6=7, 4" },
    @{ Row = 17; Col = 2; Text = "43÷2=21, 1" },
    @{ Row = 17; Col = 3; Text = "70÷5=14, 0" },
    @{ Row = 17; Col = 4; Text = "54÷2=27, 0" },
    @{ Row = 17; Col = 5; Text = "65÷9=7, 2" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}

Write-Output "done"
